$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Update the existing last row's 3rd cell text (append the finished-option clause)
$lastRow = $t.Rows.Item($t.Rows.Count)
$t.Cell($lastRow.Index, 3).Range.Text = "- Continued working on architectural design and finished “option 4 - quiz (multiple choice)” (work history on doc)"

# 2) Append a brand-new row to the end of the table for "December 29th"
$newRow = $t.Rows.Add()
$newRow.HeadingFormat = 0

$t.Cell($newRow.Index, 1).Range.Text = "December 29th"
$t.Cell($newRow.Index, 3).Range.Text = "- Continued working on architectural design and finished “option 2 - lesson” (work history on doc)"
